# Insert a new data row at row 134 (shifting existing rows 134..240 down to
# 135..241), then populate the new row with its values.
#
# This reproduces the author's edit: a new price-report entry (date 44741,
# "Perú" origin, volume 480) was inserted into the "Poroto verde" sheet right
# before the former row 134, pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 134:240 down to 135:241, freeing up row 134 for the new entry.
$ws.Rows.Item(134).Insert()

# Fill in the new row 134 with its data.
$ws.Range("A134").Value = 8
$ws.Range("B134").Value = "Terminal La Palmera de La Serena"
$ws.Range("C134").Value = "Coquimbo"
$ws.Range("D134").Value = 44741
$ws.Range("E134").Value = 4
$ws.Range("F134").Value = 100112031
$ws.Range("G134").Value = "Poroto verde"
$ws.Range("H134").Value = "Magnum"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 480
$ws.Range("K134").Value = 25000
$ws.Range("L134").Value = 26000
$ws.Range("M134").Value = 25500
$ws.Range("N134").Value = "`$/malla 25 kilos"
$ws.Range("O134").Value = "Perú"
$ws.Range("P134").Value = 1020
$ws.Range("Q134").Value = 25
$ws.Range("R134").Value = "Hortaliza"
